$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Special")
$v = $ws.Range("B3").Value()
Write-Host $v
